$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.347.10"
$ws.Range("E2").Value = "'  +0.53%  "
$ws.Range("D3").Value = "'1.876.63"
$ws.Range("E3").Value = "'  +0.90%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D5").Value = "'0.7113"
$ws.Range("E5").Value = "'  -0.45%  "
$ws.Range("D6").Value = "'242.21"
$ws.Range("E6").Value = "'  +0.73%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("D8").Value = "'0.3110"
$ws.Range("E8").Value = "'  +1.10%  "
$ws.Range("D9").Value = "'0.07749"
$ws.Range("E9").Value = "'  +0.01%  "
$ws.Range("D10").Value = "'25.07"
$ws.Range("E10").Value = "'  +0.34%  "
$ws.Range("D11").Value = "'0.08461"
$ws.Range("E11").Value = "'  +2.48%  "
$ws.Range("D12").Value = "'1.911.15"
$ws.Range("E12").Value = "'  +2.91%  "
$ws.Range("D13").Value = "'5.212"
$ws.Range("E13").Value = "'  -0.37%  "
$ws.Range("D14").Value = "'0.7116"
$ws.Range("E14").Value = "'  -0.35%  "
$ws.Range("D15").Value = "'91.40"
$ws.Range("E15").Value = "'  +1.31%  "
$ws.Range("D16").Value = "'29.358.24"
$ws.Range("E16").Value = "'  +0.63%  "
$ws.Range("D17").Value = "'0.000008282"
$ws.Range("E17").Value = "'  +6.11%  "
$ws.Range("D18").Value = "'6.009"
$ws.Range("E18").Value = "'  +2.43%  "
$ws.Range("D19").Value = "'242.41"
$ws.Range("E19").Value = "'  -0.68%  "
$ws.Range("E20").Value = "'  +0.62%  "
$ws.Range("D21").Value = "'2.129.24"
$ws.Range("E21").Value = "'  +1.41%  "
$ws.Range("D22").Value = "'0.9999"
$ws.Range("E22").Value = "'  -0.06%  "
$ws.Range("D23").Value = "'7.854"
$ws.Range("E23").Value = "'  -1.07%  "
$ws.Range("E24").Value = "'  -0.03%  "
$ws.Range("D25").Value = "'0.1606"
$ws.Range("E25").Value = "'  +1.42%  "
$ws.Range("D26").Value = "'162.63"
$ws.Range("E26").Value = "'  -0.03%  "
$ws.Range("D27").Value = "'9.018"
$ws.Range("E27").Value = "'  +1.12%  "
$ws.Range("E28").Value = "'  +1.18%  "
$ws.Range("D29").Value = "'1.512"
$ws.Range("E29").Value = "'  +1.26%  "
$ws.Range("D30").Value = "'4.407"
$ws.Range("E30").Value = "'  +0.48%  "
$ws.Range("D31").Value = "'4.335"
$ws.Range("E31").Value = "'  +4.66%  "
$ws.Range("D32").Value = "'1.277"
$ws.Range("E32").Value = "'  -2.82%  "
$ws.Range("D33").Value = "'0.05254"
$ws.Range("E33").Value = "'  +1.29%  "
$ws.Range("D34").Value = "'1.931"
$ws.Range("E34").Value = "'  +1.27%  "
$ws.Range("D35").Value = "'1.176"
$ws.Range("E35").Value = "'  +0.15%  "
$ws.Range("D36").Value = "'0.7412"
$ws.Range("E36").Value = "'  +1.83%  "
$ws.Range("E37").Value = "'  +0.18%  "
$ws.Range("D38").Value = "'0.01866"
$ws.Range("E38").Value = "'  +0.97%  "
$ws.Range("E39").Value = "'  +1.41%  "
$ws.Range("D40").Value = "'1.174.20"
$ws.Range("E40").Value = "'  +1.68%  "
$ws.Range("D41").Value = "'6.384"
$ws.Range("E41").Value = "'  +4.61%  "
$ws.Range("D42").Value = "'73.14"
$ws.Range("E42").Value = "'  +1.05%  "
$ws.Range("D43").Value = "'0.8856"
$ws.Range("E43").Value = "'  -1.99%  "
$ws.Range("D44").Value = "'106.42"
$ws.Range("E44").Value = "'  +4.66%  "
$ws.Range("D45").Value = "'0.9996"
$ws.Range("E45").Value = "'  -0.02%  "
$ws.Range("D46").Value = "'2.025.48"
$ws.Range("E46").Value = "'  +0.95%  "
$ws.Range("D47").Value = "'1.812"
$ws.Range("E47").Value = "'  +2.71%  "
$ws.Range("D48").Value = "'0.5203"
$ws.Range("E48").Value = "'  -0.64%  "
$ws.Range("E49").Value = "'  +1.44%  "
$ws.Range("E50").Value = "'  +0.68%  "
$ws.Range("E51").Value = "'  +1.09%  "
